# Apply cryptos list update (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($range, $text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Row 2
Set-TextCell $ws.Range("D2") '41.529.94'
Set-TextCell $ws.Range("E2") '  +1.55%  '

# Row 3
Set-TextCell $ws.Range("D3") '2.484.30'
Set-TextCell $ws.Range("E3") '  +1.51%  '

# Row 4
Set-TextCell $ws.Range("D4") '0.998'
Set-TextCell $ws.Range("E4") '  -0.13%  '

# Row 5
Set-TextCell $ws.Range("D5") '312.82'
Set-TextCell $ws.Range("E5") '  +0.95%  '

# Row 6
Set-TextCell $ws.Range("D6") '92.96'
Set-TextCell $ws.Range("E6") '  -0.84%  '

# Row 7
Set-TextCell $ws.Range("E7") '  -0.55%  '

# Row 8
Set-TextCell $ws.Range("E8") '  -0.20%  '

# Row 9
Set-TextCell $ws.Range("E9") '  -0.81%  '

# Row 10
Set-TextCell $ws.Range("D10") '32.67'
Set-TextCell $ws.Range("E10") '  -2.21%  '

# Row 11
Set-TextCell $ws.Range("D11") '0.0781'
Set-TextCell $ws.Range("E11") '  +0.38%  '

# Row 12
Set-TextCell $ws.Range("E12") '  +2.30%  '

# Row 13
Set-TextCell $ws.Range("D13") '2.864.39'
Set-TextCell $ws.Range("E13") '  +1.31%  '

# Row 14
Set-TextCell $ws.Range("D14") '6.85'
Set-TextCell $ws.Range("E14") '  -1.02%  '

# Row 15
Set-TextCell $ws.Range("D15") '15.45'
Set-TextCell $ws.Range("E15") '  +7.49%  '

# Row 16
Set-TextCell $ws.Range("D16") '2.476.53'
Set-TextCell $ws.Range("E16") '  +0.71%  '

# Row 17
Set-TextCell $ws.Range("E17") '  -3.63%  '

# Row 18
Set-TextCell $ws.Range("D18") '41.675.35'
Set-TextCell $ws.Range("E18") '  +1.84%  '

# Row 19
Set-TextCell $ws.Range("D19") '6.32'
Set-TextCell $ws.Range("E19") '  -0.25%  '

# Row 20
Set-TextCell $ws.Range("D20") '0.0₃0922'
Set-TextCell $ws.Range("E20") '  +1.17%  '

# Row 21
Set-TextCell $ws.Range("D21") '70.60'
Set-TextCell $ws.Range("E21") '  +5.54%  '

# Row 22
Set-TextCell $ws.Range("D22") '11.16'
Set-TextCell $ws.Range("E22") '  -3.18%  '

# Row 23
Set-TextCell $ws.Range("D23") '235.43'
Set-TextCell $ws.Range("E23") '  -0.08%  '

# Row 24
Set-TextCell $ws.Range("D24") '2.71'
Set-TextCell $ws.Range("E24") '  -1.77%  '

# Row 25
Set-TextCell $ws.Range("E25") '  +0.04%  '

# Row 26
Set-TextCell $ws.Range("E26") '  -0.58%  '

# Row 27
Set-TextCell $ws.Range("D27") '24.73'
Set-TextCell $ws.Range("E27") '  +1.29%  '

# Row 28
Set-TextCell $ws.Range("E28") '  +0.52%  '

# Row 29
Set-TextCell $ws.Range("D29") '9.64'
Set-TextCell $ws.Range("E29") '  -0.13%  '

# Row 30
Set-TextCell $ws.Range("D30") '36.22'
Set-TextCell $ws.Range("E30") '  +0.49%  '

# Row 31
Set-TextCell $ws.Range("D31") '154.32'
Set-TextCell $ws.Range("E31") '  +1.11%  '

# Row 32
Set-TextCell $ws.Range("D32") '5.41'
Set-TextCell $ws.Range("E32") '  -2.63%  '

# Row 33
Set-TextCell $ws.Range("D33") '2.57'
Set-TextCell $ws.Range("E33") '  -0.91%  '

# Row 34
Set-TextCell $ws.Range("D34") '18.15'
Set-TextCell $ws.Range("E34") '  +6.01%  '

# Row 35
Set-TextCell $ws.Range("E35") '  +0.97%  '

# Row 36
Set-TextCell $ws.Range("E36") '  -2.32%  '

# Row 37
Set-TextCell $ws.Range("E37") '  -1.59%  '

# Row 38
Set-TextCell $ws.Range("B38") 'Kaspa'
Set-TextCell $ws.Range("C38") 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextCell $ws.Range("D38") '0.105'
Set-TextCell $ws.Range("E38") '  +2.14%  '

# Row 39
Set-TextCell $ws.Range("B39") 'ARBITRUM'
Set-TextCell $ws.Range("C39") 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell $ws.Range("D39") '1.83'
Set-TextCell $ws.Range("E39") '  -3.28%  '

# Row 40
Set-TextCell $ws.Range("E40") '  -0.23%  '

# Row 41
Set-TextCell $ws.Range("D41") '4.11'
Set-TextCell $ws.Range("E41") '  -0.30%  '

# Row 42
Set-TextCell $ws.Range("E42") '  -0.13%  '

# Row 43
Set-TextCell $ws.Range("D43") '19.76'
Set-TextCell $ws.Range("E43") '  -5.86%  '

# Row 44
Set-TextCell $ws.Range("D44") '1.955.86'
Set-TextCell $ws.Range("E44") '  -0.06%  '

# Row 45
Set-TextCell $ws.Range("D45") '0.0285'
Set-TextCell $ws.Range("E45") '  +0.55%  '

# Row 46
Set-TextCell $ws.Range("E46") '  -2.43%  '

# Row 47
Set-TextCell $ws.Range("E47") '  +1.83%  '

# Row 48
Set-TextCell $ws.Range("D48") '2.726.11'
Set-TextCell $ws.Range("E48") '  +1.30%  '

# Row 49
Set-TextCell $ws.Range("D49") '96.16'
Set-TextCell $ws.Range("E49") '  -0.89%  '

# Row 50
Set-TextCell $ws.Range("D50") '67.08'
Set-TextCell $ws.Range("E50") '  -2.83%  '

# Row 51
Set-TextCell $ws.Range("D51") '73.14'
Set-TextCell $ws.Range("E51") '  -3.91%  '
